$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 95, shifting existing rows 95-105 down to 97-107
$ws.Range("A95:T96").EntireRow.Insert()

# Row 95: new data row, date 44449, Especial quality
$ws.Range("A95").Value = 4
$ws.Range("B95").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C95").Value = "Los Lagos"
$ws.Range("D95").Value = 44449
$ws.Range("D95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E95").Value = 10
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100101
$ws.Range("H95").Value = "Berries"
$ws.Range("I95").Value = 100101007
$ws.Range("J95").Value = "Kiwi"
$ws.Range("K95").Value = "Hayward"
$ws.Range("L95").Value = "Especial"
$ws.Range("M95").Value = 300
$ws.Range("N95").Value = 20000
$ws.Range("O95").Value = 20000
$ws.Range("P95").Value = 20000
$ws.Range("Q95").Value = "$/caja 15 kilos"
$ws.Range("R95").Value = "Provincia de Curicó"
$ws.Range("S95").Value = 1333
$ws.Range("T95").Value = 15

# Row 96: new data row, date 44449, Primera quality
$ws.Range("A96").Value = 4
$ws.Range("B96").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C96").Value = "Los Lagos"
$ws.Range("D96").Value = 44449
$ws.Range("D96").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E96").Value = 10
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100101
$ws.Range("H96").Value = "Berries"
$ws.Range("I96").Value = 100101007
$ws.Range("J96").Value = "Kiwi"
$ws.Range("K96").Value = "Hayward"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 150
$ws.Range("N96").Value = 14000
$ws.Range("O96").Value = 14000
$ws.Range("P96").Value = 14000
$ws.Range("Q96").Value = "$/caja 15 kilos"
$ws.Range("R96").Value = "Provincia de Curicó"
$ws.Range("S96").Value = 933
$ws.Range("T96").Value = 15
